$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1781.875
$ws.Range("I101").Value = 1926.7142
$ws.Range("J101").Value = 768
$ws.Range("K101").Value = 5780.142599999999
$ws.Range("L101").Value = 2304
$ws.Range("M101").Value = -4158.142599999999
$ws.Range("N101").Value = -5548
$ws.Range("H125").Value = 125001400
$ws.Range("I125").Value = 500000320
$ws.Range("K125").Value = 4500002880
$ws.Range("M125").Value = -4500000420
$ws.Range("H129").Value = 13691.09
$ws.Range("I129").Value = 501.35294
$ws.Range("J129").Value = 17366.918
$ws.Range("K129").Value = 1504.05882
$ws.Range("L129").Value = 52100.754
$ws.Range("M129").Value = 3495.94118
$ws.Range("N129").Value = -62100.754
$ws.Range("H132").Value = 3573005.8
$ws.Range("I132").Value = 4763614.5
$ws.Range("J132").Value = 1180.1
$ws.Range("K132").Value = 14290843.5
$ws.Range("L132").Value = 3540.3
$ws.Range("M132").Value = -14288313.5
$ws.Range("N132").Value = -8600.299999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14149.719
$ws.Range("I32").Value = 15427.506
$ws.Range("K32").Value = 15427.506
$ws.Range("M32").Value = -15140.506
$ws.Range("H45").Value = 1144.0625
$ws.Range("I45").Value = 1019.63635
$ws.Range("J45").Value = 1417.8
$ws.Range("K45").Value = 1019.63635
$ws.Range("L45").Value = 1417.8
$ws.Range("M45").Value = -642.63635
$ws.Range("N45").Value = -2171.8
$ws.Range("H61").Value = 2488.25
$ws.Range("I61").Value = 1712.6666
$ws.Range("J61").Value = 3780.889
$ws.Range("K61").Value = 1712.6666
$ws.Range("L61").Value = 3780.889
$ws.Range("M61").Value = -1500.6666
$ws.Range("N61").Value = -4204.889
$ws.Range("H74").Value = 746.80646
$ws.Range("I74").Value = 659.6539
$ws.Range("K74").Value = 659.6539
$ws.Range("M74").Value = 214.3461
$ws.Range("H75").Value = 156000
$ws.Range("I75").Value = 12000
$ws.Range("K75").Value = 12000
$ws.Range("M75").Value = -11126
$ws.Range("H76").Value = 25290
$ws.Range("J76").Value = 25290
$ws.Range("L76").Value = 25290
$ws.Range("N76").Value = -25966
$ws.Range("H77").Value = 746.80646
$ws.Range("I77").Value = 659.6539
$ws.Range("K77").Value = 3298.2695
$ws.Range("M77").Value = 1069.7305
$ws.Range("H78").Value = 156000
$ws.Range("I78").Value = 12000
$ws.Range("K78").Value = 36000
$ws.Range("M78").Value = -31632
$ws.Range("H79").Value = 25290
$ws.Range("J79").Value = 25290
$ws.Range("L79").Value = 25290
$ws.Range("N79").Value = -27630
$ws.Range("H97").Value = 721.375
$ws.Range("I97").Value = 878
$ws.Range("J97").Value = 460.33334
$ws.Range("K97").Value = 878
$ws.Range("L97").Value = 460.33334
$ws.Range("M97").Value = -382
$ws.Range("N97").Value = -1452.33334
$ws.Range("H136").Value = 2488.25
$ws.Range("I136").Value = 1712.6666
$ws.Range("J136").Value = 3780.889
$ws.Range("K136").Value = 5137.9998
$ws.Range("L136").Value = 11342.667
$ws.Range("M136").Value = -2587.9998
$ws.Range("N136").Value = -16442.667

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 866.5
$ws.Range("I94").Value = 1450
$ws.Range("J94").Value = 574.75
$ws.Range("K94").Value = 1450
$ws.Range("L94").Value = 574.75
$ws.Range("M94").Value = -999
$ws.Range("N94").Value = -1476.75
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H134").Value = 25952.215
$ws.Range("I134").Value = 35155.965
$ws.Range("J134").Value = 2942.8333
$ws.Range("K134").Value = 105467.895
$ws.Range("L134").Value = 8828.499899999999
$ws.Range("M134").Value = -102932.895
$ws.Range("N134").Value = -13898.4999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H76").Value = 500
$ws.Range("I76").Value = 500
$ws.Range("K76").Value = 500
$ws.Range("M76").Value = -185
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0
$ws.Range("H79").Value = 500
$ws.Range("I79").Value = 500
$ws.Range("K79").Value = 500
$ws.Range("M79").Value = 592
$ws.Range("H132").Value = 1621.75
$ws.Range("I132").Value = 868.2857
$ws.Range("J132").Value = 2676.6
$ws.Range("K132").Value = 2604.8571
$ws.Range("L132").Value = 8029.799999999999
$ws.Range("M132").Value = -74.85710000000017
$ws.Range("N132").Value = -13089.8
$ws.Range("H134").Value = 1148.4445
$ws.Range("I134").Value = 1075.7407
$ws.Range("J134").Value = 1366.5555
$ws.Range("K134").Value = 3227.2221
$ws.Range("L134").Value = 4099.666499999999
$ws.Range("M134").Value = -692.2221
$ws.Range("N134").Value = -9169.666499999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 600
$ws.Range("I109").Value = 600
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 1800
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").Value = -760
$ws.Range("H113").Value = 527.8222
$ws.Range("I113").Value = 524
$ws.Range("J113").Value = 533.55554
$ws.Range("K113").Value = 1572
$ws.Range("L113").Value = 1600.66662
$ws.Range("M113").Value = 598
$ws.Range("N113").Value = -5940.66662
$ws.Range("H131").Value = 2852552.8
$ws.Range("J131").Value = 4116133.5
$ws.Range("L131").Value = 12348400.5
$ws.Range("N131").Value = -12358480.5
$ws.Range("H132").Value = 2335.2942
$ws.Range("I132").Value = 1182.8572
$ws.Range("K132").Value = 10645.7148
$ws.Range("M132").Value = -8115.7148

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4780
$ws.Range("I80").Value = 2760
$ws.Range("J80").Value = 6800
$ws.Range("K80").Value = 2760
$ws.Range("L80").Value = 6800
$ws.Range("M80").Value = -1762
$ws.Range("N80").Value = -8796
$ws.Range("H83").Value = 4780
$ws.Range("I83").Value = 2760
$ws.Range("J83").Value = 6800
$ws.Range("K83").Value = 13800
$ws.Range("L83").Value = 34000
$ws.Range("M83").Value = -8808
$ws.Range("N83").Value = -43984
$ws.Range("H97").Value = 981.7646999999999
$ws.Range("I97").Value = 860.7692
$ws.Range("J97").Value = 1375
$ws.Range("K97").Value = 860.7692
$ws.Range("L97").Value = 1375
$ws.Range("M97").Value = -364.7692
$ws.Range("N97").Value = -2367
$ws.Range("H104").Value = 29636.625
$ws.Range("I104").Value = 23333
$ws.Range("J104").Value = 30537.143
$ws.Range("K104").Value = 23333
$ws.Range("L104").Value = 30537.143
$ws.Range("M104").Value = -19839
$ws.Range("N104").Value = -37525.143
$ws.Range("H132").Value = 57888
$ws.Range("I132").Value = 84756
$ws.Range("J132").Value = 4152
$ws.Range("K132").Value = 254268
$ws.Range("L132").Value = 12456
$ws.Range("M132").Value = -251738
$ws.Range("N132").Value = -17516

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H82").Value = 2625
$ws.Range("I82").Value = 2625
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2625
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2264
$ws.Range("H85").Value = 2625
$ws.Range("I85").Value = 2625
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2625
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -1377
$ws.Range("H132").Value = 1574.0638
$ws.Range("I132").Value = 1348.3529
$ws.Range("J132").Value = 2164.3845
$ws.Range("K132").Value = 4045.0587
$ws.Range("L132").Value = 6493.1535
$ws.Range("M132").Value = -1515.0587
$ws.Range("N132").Value = -11553.1535
$ws.Range("H133").Value = 23756.5
$ws.Range("J133").Value = 23756.5
$ws.Range("L133").Value = 23756.5
$ws.Range("N133").Value = -28816.5
$ws.Range("H136").Value = 7260.6113
$ws.Range("I136").Value = 9057.583000000001
$ws.Range("K136").Value = 27172.749
$ws.Range("M136").Value = -24622.749

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 35714932
$ws.Range("I96").Value = 62500508
$ws.Range("K96").Value = 62500508
$ws.Range("M96").Value = -62499135
$ws.Range("H107").Value = 458.33334
$ws.Range("I107").Value = 414.4
$ws.Range("J107").Value = 513.25
$ws.Range("K107").Value = 1243.2
$ws.Range("L107").Value = 1539.75
$ws.Range("M107").Value = 676.8000000000002
$ws.Range("N107").Value = -5379.75
$ws.Range("H132").Value = 1613.2559
$ws.Range("I132").Value = 1416.9688
$ws.Range("J132").Value = 2184.2727
$ws.Range("K132").Value = 4250.9064
$ws.Range("L132").Value = 6552.8181
$ws.Range("M132").Value = -1720.9064
$ws.Range("N132").Value = -11612.8181
$ws.Range("H136").Value = 5188.5757
$ws.Range("I136").Value = 5723.9614
$ws.Range("J136").Value = 3200
$ws.Range("K136").Value = 17171.8842
$ws.Range("L136").Value = 9600
$ws.Range("M136").Value = -14621.8842
$ws.Range("N136").Value = -14700

Write-Host "Applied all Garuda_Profits updates"